# Apply price/volume updates to cryptos worksheet, per commit:
# "Updated cryptos list on Tue Sep  5 14:08:35 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '25.834.10'
$ws.Range("E2").Value = '  -0.25%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.635.84'
$ws.Range("E3").Value = '  +0.04%  '

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '215.19'
$ws.Range("E5").Value = '  -0.54%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.5038'
$ws.Range("E6").Value = '  -0.67%  '

# Row 7
$ws.Range("E7").Value = '  +0.27%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.2571'
$ws.Range("E8").Value = '  -0.29%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06412'
$ws.Range("E9").Value = '  +0.79%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '19.92'
$ws.Range("E10").Value = '  +1.66%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07797'
$ws.Range("E11").Value = '  +0.53%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '4.280'
$ws.Range("E12").Value = '  +0.47%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.863.00'
$ws.Range("E13").Value = '  +0.10%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '1.633.74'
$ws.Range("E14").Value = '  -0.03%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.5592'
$ws.Range("E15").Value = '  +1.35%  '

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.0₅7616'
$ws.Range("E16").Value = '  -1.23%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '62.91'
$ws.Range("E17").Value = '  -1.74%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '25.859.28'
$ws.Range("E18").Value = '  -0.23%  '

# Row 19
$ws.Range("E19").Value = '  +0.13%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '194.36'
$ws.Range("E20").Value = '  -0.03%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '4.324'
$ws.Range("E21").Value = '  -2.78%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.869'
$ws.Range("E22").Value = '  -0.27%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.090'
$ws.Range("E23").Value = '  +0.46%  '

# Row 24
$ws.Range("E24").Value = '  +0.17%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '1.779'
$ws.Range("E25").Value = '  -6.77%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '140.24'
$ws.Range("E26").Value = '  -1.70%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.1252'
$ws.Range("E27").Value = '  +1.32%  '

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '6.815'
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '15.45'
$ws.Range("E29").Value = '  -0.82%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.240'
$ws.Range("E30").Value = '  -0.38%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.04884'
$ws.Range("E31").Value = '  +0.36%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.293'
$ws.Range("E32").Value = '  +1.29%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.218'
$ws.Range("E33").Value = '  +0.91%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.565'
$ws.Range("E34").Value = '  +1.42%  '

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.381'
$ws.Range("E35").Value = '  +0.30%  '

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.9002'
$ws.Range("E36").Value = '  -0.56%  '

# Row 37
$ws.Range("E37").Value = '  +0.29%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.5525'
$ws.Range("E38").Value = '  +0.57%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.125.00'
$ws.Range("E39").Value = '  +0.20%  '

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.01560'
$ws.Range("E40").Value = '  +0.06%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.9966'
$ws.Range("E41").Value = '  -0.49%  '

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '5.528'
$ws.Range("E42").Value = '  -0.90%  '

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '97.95'
$ws.Range("E44").Value = '  +0.51%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.773.45'
$ws.Range("E45").Value = '  -0.03%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0₈113'
$ws.Range("E46").Value = '  -6.51%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '55.39'
$ws.Range("E47").Value = '  +0.91%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.4264'
$ws.Range("E48").Value = '  -4.36%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '7.721'
$ws.Range("E49").Value = '  +2.44%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.05036'
$ws.Range("E50").Value = '  -2.25%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$ws.Range("E51").Value = '  +0.34%  '
